$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reordered EC (estado de cuenta) data: newest employee/period records first.
# Column B=TipoDoc, C=NumDoc, D=Nombre, E=Periodo, F=ValorMora, G=SalarioBasico
$data = @(
    ,@("CC", "1047365927", "JOHORMAN GONZALEZ GOMEZ", "2311", 61290, 1820000)
    ,@("CC", "1047365927", "JOHORMAN GONZALEZ GOMEZ", "2310", 72800, 1820000)
    ,@("CC", "1047365927", "JOHORMAN GONZALEZ GOMEZ", "2309", 72800, 1820000)
    ,@("CC", "1047365927", "JOHORMAN GONZALEZ GOMEZ", "2308", 72800, 1820000)
    ,@("CC", "1047365927", "JOHORMAN GONZALEZ GOMEZ", "2307", 72800, 1820000)
    ,@("CC", "1047365927", "JOHORMAN GONZALEZ GOMEZ", "2306", 72800, 1820000)
    ,@("CC", "1047365927", "JOHORMAN GONZALEZ GOMEZ", "2305", 72800, 1820000)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2311", 16386, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2310", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2309", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2308", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2307", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2306", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2305", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2304", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2303", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2302", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2301", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2212", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2211", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2210", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2209", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2208", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2207", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2206", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2205", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2204", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2203", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2202", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2201", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2112", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2111", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2110", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2109", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2108", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2107", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2106", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2105", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2104", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2103", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2102", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2101", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2012", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2011", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2010", 17556, 438901)
    ,@("CC", "1067866451", "CAMILO ANDRES KERGUELEN MANJARREZ", "2009", 17556, 438901)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}
